$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01692966666666667
$ws.Range("H2").Value = 0.050789
$ws.Range("I2").Value = 0.007576566749688815
$ws.Range("J2").Value = 0.008018731266168228
$ws.Range("M2").Value = 1.521496
$ws.Range("N2").Value = 4.564488
$ws.Range("O2").Value = 0.03673310907796454
$ws.Range("P2").Value = 0.04709307785339354
$ws.Range("Q2").Value = 0.02575842011466667
$ws.Range("R2").Value = 0.231825781032
$ws.Range("S2").Value = 0.0002783108528527986
$ws.Range("T2").Value = 0.0003776267358031013

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01692966666666667
$ws.Range("H3").Value = 0.050789
$ws.Range("I3").Value = 0.007576566749688815
$ws.Range("J3").Value = 0.008018731266168228
$ws.Range("O3").Value = 0.2878223347346259
$ws.Range("P3").Value = 0.3689978865887579
$ws.Range("Q3").Value = 0.2018301418685556
$ws.Range("R3").Value = 1.816471276817
$ws.Range("S3").Value = 0.002180705131168171
$ws.Range("T3").Value = 0.002958894890339271

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01692966666666667
$ws.Range("H4").Value = 0.050789
$ws.Range("I4").Value = 0.007576566749688815
$ws.Range("J4").Value = 0.008018731266168228
$ws.Range("M4").Value = 0.108844
$ws.Range("N4").Value = 0.326532
$ws.Range("O4").Value = 0.002627794305395461
$ws.Range("P4").Value = 0.003368920434805459
$ws.Range("Q4").Value = 0.001842692638666667
$ws.Range("R4").Value = 0.016584233748
$ws.Range("S4").Value = 0.00001990965895928087
$ws.Range("T4").Value = 0.0000270144676238076

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01692966666666667
$ws.Range("H5").Value = 0.050789
$ws.Range("I5").Value = 0.007576566749688815
$ws.Range("J5").Value = 0.008018731266168228
$ws.Range("M5").Value = 27.3360495
$ws.Range("N5").Value = 54.672099
$ws.Range("O5").Value = 0.659967616112128
$ws.Range("P5").Value = 0.5640670792902598
$ws.Range("Q5").Value = 0.4627902060185
$ws.Range("R5").Value = 2.776741236111
$ws.Range("S5").Value = 0.005000288696106541
$ws.Range("T5").Value = 0.004523102324920999

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01692966666666667
$ws.Range("H6").Value = 0.050789
$ws.Range("I6").Value = 0.007576566749688815
$ws.Range("J6").Value = 0.008018731266168228
$ws.Range("M6").Value = 0.5322153333333334
$ws.Range("N6").Value = 1.596646
$ws.Range("O6").Value = 0.01284914576988608
$ws.Range("P6").Value = 0.0164730358327833
$ws.Range("Q6").Value = 0.009010228188222223
$ws.Range("R6").Value = 0.08109205369400001
$ws.Range("S6").Value = 0.00009735241060202359
$ws.Range("T6").Value = 0.000132092847481049

$ws.Range("G7").Value = 0.9705113333333335
$ws.Range("H7").Value = 2.911534000000001
$ws.Range("I7").Value = 0.434334830277983
$ws.Range("J7").Value = 0.459682386310261
$ws.Range("M7").Value = 1.521496
$ws.Range("N7").Value = 4.564488
$ws.Range("O7").Value = 0.03673310907796454
$ws.Range("P7").Value = 0.04709307785339354
$ws.Range("Q7").Value = 1.476629111621333
$ws.Range("R7").Value = 13.289662004592
$ws.Range("S7").Value = 0.01595446869696037
$ws.Range("T7").Value = 0.02164785840634285

$ws.Range("G8").Value = 0.9705113333333335
$ws.Range("H8").Value = 2.911534000000001
$ws.Range("I8").Value = 0.434334830277983
$ws.Range("J8").Value = 0.459682386310261
$ws.Range("O8").Value = 0.2878223347346259
$ws.Range("P8").Value = 0.3689978865887579
$ws.Range("Q8").Value = 11.57012975792245
$ws.Range("S8").Value = 0.1250112649071766
$ws.Range("T8").Value = 0.1696218290505633

$ws.Range("G9").Value = 0.9705113333333335
$ws.Range("H9").Value = 2.911534000000001
$ws.Range("I9").Value = 0.434334830277983
$ws.Range("J9").Value = 0.459682386310261
$ws.Range("M9").Value = 0.108844
$ws.Range("N9").Value = 0.326532
$ws.Range("O9").Value = 0.002627794305395461
$ws.Range("P9").Value = 0.003368920434805459
$ws.Range("Q9").Value = 0.1056343355653333
$ws.Range("R9").Value = 0.9507090200880002
$ws.Range("S9").Value = 0.001141342593639388
$ws.Range("T9").Value = 0.001548633384760776

$ws.Range("G10").Value = 0.9705113333333335
$ws.Range("H10").Value = 2.911534000000001
$ws.Range("I10").Value = 0.434334830277983
$ws.Range("J10").Value = 0.459682386310261
$ws.Range("M10").Value = 27.3360495
$ws.Range("N10").Value = 54.672099
$ws.Range("O10").Value = 0.659967616112128
$ws.Range("P10").Value = 0.5640670792902598
$ws.Range("Q10").Value = 26.52994584831101
$ws.Range("R10").Value = 159.179675089866
$ws.Range("S10").Value = 0.2866469225330261
$ws.Range("T10").Value = 0.2592917010472058

$ws.Range("G11").Value = 0.9705113333333335
$ws.Range("H11").Value = 2.911534000000001
$ws.Range("I11").Value = 0.434334830277983
$ws.Range("J11").Value = 0.459682386310261
$ws.Range("M11").Value = 0.5322153333333334
$ws.Range("N11").Value = 1.596646
$ws.Range("O11").Value = 0.01284914576988608
$ws.Range("P11").Value = 0.0164730358327833
$ws.Range("Q11").Value = 0.5165210127737779
$ws.Range("R11").Value = 4.648689114964001
$ws.Range("S11").Value = 0.005580831547180536
$ws.Range("T11").Value = 0.007572364421388266

$ws.Range("G12").Value = 0.5101613333333334
$ws.Range("H12").Value = 1.530484
$ws.Range("I12").Value = 0.2283134967282431
$ws.Range("J12").Value = 0.2416377543005417
$ws.Range("M12").Value = 1.521496
$ws.Range("N12").Value = 4.564488
$ws.Range("O12").Value = 0.03673310907796454
$ws.Range("P12").Value = 0.04709307785339354
$ws.Range("Q12").Value = 0.7762084280213334
$ws.Range("R12").Value = 6.985875852192
$ws.Range("S12").Value = 0.008386664579290054
$ws.Range("T12").Value = 0.01137946557559459

$ws.Range("G13").Value = 0.5101613333333334
$ws.Range("H13").Value = 1.530484
$ws.Range("I13").Value = 0.2283134967282431
$ws.Range("J13").Value = 0.2416377543005417
$ws.Range("O13").Value = 0.2878223347346259
$ws.Range("P13").Value = 0.3689978865887579
$ws.Range("Q13").Value = 6.081982375072445
$ws.Range("R13").Value = 54.737841375652
$ws.Range("S13").Value = 0.0657137236797493
$ws.Range("T13").Value = 0.08916382065695344

$ws.Range("G14").Value = 0.5101613333333334
$ws.Range("H14").Value = 1.530484
$ws.Range("I14").Value = 0.2283134967282431
$ws.Range("J14").Value = 0.2416377543005417
$ws.Range("M14").Value = 0.108844
$ws.Range("N14").Value = 0.326532
$ws.Range("O14").Value = 0.002627794305395461
$ws.Range("P14").Value = 0.003368920434805459
$ws.Range("Q14").Value = 0.05552800016533333
$ws.Range("R14").Value = 0.499752001488
$ws.Range("S14").Value = 0.0005999609065474024
$ws.Range("T14").Value = 0.0008140583682835957

$ws.Range("G15").Value = 0.5101613333333334
$ws.Range("H15").Value = 1.530484
$ws.Range("I15").Value = 0.2283134967282431
$ws.Range("J15").Value = 0.2416377543005417
$ws.Range("M15").Value = 27.3360495
$ws.Range("N15").Value = 54.672099
$ws.Range("O15").Value = 0.659967616112128
$ws.Range("P15").Value = 0.5640670792902598
$ws.Range("Q15").Value = 13.945795460986
$ws.Range("R15").Value = 83.674772765916
$ws.Range("S15").Value = 0.1506795141619627
$ws.Range("T15").Value = 0.136299902314564

$ws.Range("G16").Value = 0.5101613333333334
$ws.Range("H16").Value = 1.530484
$ws.Range("I16").Value = 0.2283134967282431
$ws.Range("J16").Value = 0.2416377543005417
$ws.Range("M16").Value = 0.5322153333333334
$ws.Range("N16").Value = 1.596646
$ws.Range("O16").Value = 0.01284914576988608
$ws.Range("P16").Value = 0.0164730358327833
$ws.Range("Q16").Value = 0.2715156840737778
$ws.Range("R16").Value = 2.443641156664
$ws.Range("S16").Value = 0.002933633400693605
$ws.Range("T16").Value = 0.003980507385146111

$ws.Range("G17").Value = 0.369637
$ws.Range("H17").Value = 0.739274
$ws.Range("I17").Value = 0.1654243676970244
$ws.Range("J17").Value = 0.1167189654859369
$ws.Range("M17").Value = 1.521496
$ws.Range("N17").Value = 4.564488
$ws.Range("O17").Value = 0.03673310907796454
$ws.Range("P17").Value = 0.04709307785339354
$ws.Range("Q17").Value = 0.562401216952
$ws.Range("R17").Value = 3.374407301712
$ws.Range("S17").Value = 0.006076551342768113
$ws.Range("T17").Value = 0.005496655328596781

$ws.Range("G18").Value = 0.369637
$ws.Range("H18").Value = 0.739274
$ws.Range("I18").Value = 0.1654243676970244
$ws.Range("J18").Value = 0.1167189654859369
$ws.Range("O18").Value = 0.2878223347346259
$ws.Range("P18").Value = 0.3689978865887579
$ws.Range("Q18").Value = 4.406695631920333
$ws.Range("R18").Value = 26.440173791522
$ws.Range("S18").Value = 0.0476128277325568
$ws.Range("T18").Value = 0.0430690515891369

$ws.Range("G19").Value = 0.369637
$ws.Range("H19").Value = 0.739274
$ws.Range("I19").Value = 0.1654243676970244
$ws.Range("J19").Value = 0.1167189654859369
$ws.Range("M19").Value = 0.108844
$ws.Range("N19").Value = 0.326532
$ws.Range("O19").Value = 0.002627794305395461
$ws.Range("P19").Value = 0.003368920434805459
$ws.Range("Q19").Value = 0.040232769628
$ws.Range("R19").Value = 0.241396617768
$ws.Range("S19").Value = 0.0004347012114078857
$ws.Range("T19").Value = 0.000393216907954926

$ws.Range("G20").Value = 0.369637
$ws.Range("H20").Value = 0.739274
$ws.Range("I20").Value = 0.1654243676970244
$ws.Range("J20").Value = 0.1167189654859369
$ws.Range("M20").Value = 27.3360495
$ws.Range("N20").Value = 54.672099
$ws.Range("O20").Value = 0.659967616112128
$ws.Range("P20").Value = 0.5640670792902598
$ws.Range("Q20").Value = 10.1044153290315
$ws.Range("R20").Value = 40.417661316126
$ws.Range("S20").Value = 0.1091747255958613
$ws.Range("T20").Value = 0.06583732595943308

$ws.Range("G21").Value = 0.369637
$ws.Range("H21").Value = 0.739274
$ws.Range("I21").Value = 0.1654243676970244
$ws.Range("J21").Value = 0.1167189654859369
$ws.Range("M21").Value = 0.5322153333333334
$ws.Range("N21").Value = 1.596646
$ws.Range("O21").Value = 0.01284914576988608
$ws.Range("P21").Value = 0.0164730358327833
$ws.Range("Q21").Value = 0.1967264791673334
$ws.Range("R21").Value = 1.180358875004
$ws.Range("S21").Value = 0.002125561814430302
$ws.Range("T21").Value = 0.001922715700815236

$ws.Range("G22").Value = 0.367238
$ws.Range("H22").Value = 1.101714
$ws.Range("I22").Value = 0.1643507385470607
$ws.Range("J22").Value = 0.173942162637092
$ws.Range("M22").Value = 1.521496
$ws.Range("N22").Value = 4.564488
$ws.Range("O22").Value = 0.03673310907796454
$ws.Range("P22").Value = 0.04709307785339354
$ws.Range("Q22").Value = 0.558751148048
$ws.Range("R22").Value = 5.028760332432
$ws.Range("S22").Value = 0.006037113606093211
$ws.Range("T22").Value = 0.008191471807056214

$ws.Range("G23").Value = 0.367238
$ws.Range("H23").Value = 1.101714
$ws.Range("I23").Value = 0.1643507385470607
$ws.Range("J23").Value = 0.173942162637092
$ws.Range("O23").Value = 0.2878223347346259
$ws.Range("P23").Value = 0.3689978865887579
$ws.Range("Q23").Value = 4.378095511204667
$ws.Range("R23").Value = 39.40285960084201
$ws.Range("S23").Value = 0.04730381328397509
$ws.Range("T23").Value = 0.06418429040176495

$ws.Range("G24").Value = 0.367238
$ws.Range("H24").Value = 1.101714
$ws.Range("I24").Value = 0.1643507385470607
$ws.Range("J24").Value = 0.173942162637092
$ws.Range("M24").Value = 0.108844
$ws.Range("N24").Value = 0.326532
$ws.Range("O24").Value = 0.002627794305395461
$ws.Range("P24").Value = 0.003368920434805459
$ws.Range("Q24").Value = 0.039971652872
$ws.Range("R24").Value = 0.359744875848
$ws.Range("S24").Value = 0.0004318799348415044
$ws.Range("T24").Value = 0.0005859973061823538

$ws.Range("G25").Value = 0.367238
$ws.Range("H25").Value = 1.101714
$ws.Range("I25").Value = 0.1643507385470607
$ws.Range("J25").Value = 0.173942162637092
$ws.Range("M25").Value = 27.3360495
$ws.Range("N25").Value = 54.672099
$ws.Range("O25").Value = 0.659967616112128
$ws.Range("P25").Value = 0.5640670792902598
$ws.Range("Q25").Value = 10.038836146281
$ws.Range("R25").Value = 60.233016877686
$ws.Range("S25").Value = 0.1084661651251713
$ws.Range("T25").Value = 0.09811504764413582

$ws.Range("G26").Value = 0.367238
$ws.Range("H26").Value = 1.101714
$ws.Range("I26").Value = 0.1643507385470607
$ws.Range("J26").Value = 0.173942162637092
$ws.Range("M26").Value = 0.5322153333333334
$ws.Range("N26").Value = 1.596646
$ws.Range("O26").Value = 0.01284914576988608
$ws.Range("P26").Value = 0.0164730358327833
$ws.Range("Q26").Value = 0.1954496945826667
$ws.Range("R26").Value = 1.759047251244
$ws.Range("S26").Value = 0.002111766596979618
$ws.Range("T26").Value = 0.002865355477952637
